$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-12 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-13 Sunday", 2)

$d.Content.Find.Execute("422÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "230÷6=", 2)
$d.Content.Find.Execute("100÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "391÷9=", 2)
$d.Content.Find.Execute("112÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "101÷2=", 2)
$d.Content.Find.Execute("311÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "572÷5=", 2)
$d.Content.Find.Execute("719÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "761÷9=", 2)

$d.Content.Find.Execute("838÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "341÷6=", 2)
$d.Content.Find.Execute("373÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "758÷8=", 2)
$d.Content.Find.Execute("853÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "489÷3=", 2)
$d.Content.Find.Execute("606÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "546÷6=", 2)
$d.Content.Find.Execute("266÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "367÷3=", 2)

$d.Content.Find.Execute("775÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "768÷3=", 2)
$d.Content.Find.Execute("128÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "651÷5=", 2)
$d.Content.Find.Execute("384÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "596÷9=", 2)
$d.Content.Find.Execute("753÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "515÷4=", 2)
$d.Content.Find.Execute("764÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "943÷6=", 2)

$d.Content.Find.Execute("979÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "360÷2=", 2)
$d.Content.Find.Execute("320÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "573÷8=", 2)
$d.Content.Find.Execute("754÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "601÷9=", 2)
$d.Content.Find.Execute("272÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "874÷3=", 2)
$d.Content.Find.Execute("547÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "883÷3=", 2)

$d.Content.Find.Execute("796÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "221÷6=", 2)
$d.Content.Find.Execute("458÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "991÷5=", 2)
$d.Content.Find.Execute("520÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "984÷6=", 2)
$d.Content.Find.Execute("875÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "451÷9=", 2)
$d.Content.Find.Execute("583÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "102÷4=", 2)
